# Update cryptos list values to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.554.15"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.497.76"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'605.52"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'151.74"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "3.495.64"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("E11").Value = "  +7.33%  "
$ws.Range("D12").Value = "'0.431"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "'0.0000217"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "'32.37"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "4.082.48"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.611.15"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.495.37"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'6.53"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "'15.57"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").Value = "'9.84"
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").Value = "'446.09"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'0.631"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'78.04"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "3.632.31"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +5.95%  "
$ws.Range("D32").Value = "'0.169"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'25.66"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "'6.15"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'1.86"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "3.482.86"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "'7.99"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  +7.01%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'175.30"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").Value = "'0.0890"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'5.45"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'0.889"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "'30.02"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").Value = "'46.18"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "'1.30"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("D49").Value = "'2.52"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'7.62"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'0.992"
$ws.Range("E51").Value = "  -0.85%  "
